# Boletin Epi Pereira - revision_IRAhospital_semanal
# "semana 50 de 2025": append the week-50 column (BA) to the weekly
# hospital-reporting table, plus two late-arriving week-49 (AZ) values
# for institutions whose report came in after the previous publish.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: new week-number column label in row 1 (same position/style as
# the existing "1".."49" week headers in D1:AZ1).
$ws.Range("BA1").Value = 50

# Week 50 counts per reporting institution (row = institution, sparse:
# only institutions that had already reported through week 49 get a
# week-50 cell here).
$week50 = @{
    2  = 0
    3  = 0
    5  = 0
    6  = 31
    7  = 21
    8  = 12
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    19 = 0
    23 = 0
    25 = 3
    28 = 9
    29 = 0
    30 = 0
    31 = 0
    35 = 8
    36 = 0
    38 = 0
    41 = 0
    42 = 0
    43 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 0
}

foreach ($row in $week50.Keys) {
    $ws.Range("BA$row").Value = $week50[$row]
}

# Two rows (30 and 53) also received their previously-missing week-49
# (AZ) value in this same update.
$ws.Range("AZ30").Value = 2
$ws.Range("AZ53").Value = 0
